$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "negative scenario" row (row 3): a name and an email address,
# mirroring the existing row 2 pattern (name in column A, hyperlinked
# email address in column B).
$ws.Range("A3").Value = "XYZ"
$ws.Range("B3").Value = "XYZ@gmail.com"

# Turn the email in B3 into a mailto: hyperlink, same as B2.
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:XYZ@gmail.com") | Out-Null

# Move/record the active selection on the new cell.
$ws.Range("B3").Select() | Out-Null
